# The "Förändrad" (changed) date in column C for every data row (2-236)
# is updated from serial date 45190 (2023-09-21) to 45192 (2023-09-23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C236").Value = 45192
